$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Answer column updated to "Done" and gets the same highlight fill as I8/I9 ---
$ws.Range("I7").Value = "Done"
$ws.Range("I7").Interior.Color = $ws.Range("I8").Interior.Color

# --- Row 11: Note text got shorter, and row height shrinks accordingly ---
$ws.Range("I11").Value = "Pending,`n "
$ws.Rows(11).RowHeight = 28.8

# --- Row 13: Answer now "Pending", Note text updated (figma instead of Xara Web Designer) ---
$ws.Range("I13").Value = "Pending"
$ws.Range("J13").Value = "Sử dụng figma, `n  tham khảo leetcode,`n các trang web trên mạng"

# --- Row 14: new task #13 "Thiet ke database" ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Task"
$ws.Range("C14").Value = "BE"
$ws.Range("D14").Value = "Thiết kế database"
$ws.Range("E14").Value = "HoangAncient"
$ws.Range("F14").Value = 45080

# --- Row 15: new task #14 stub (only No. and Type filled so far) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Task"

# --- Selection moves to N7 ---
$null = $ws.Range("N7").Select()
